# Update the "Förändrad" (Changed) date column (C) from 2023-09-01 (45170)
# to 2023-09-05 (45174) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp from bottom of column B

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value2 = 45174
    }
}
